$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the text of the existing last log entry (C111) ---
$ws.Range("C111").Value = "Reset, Base, Header, LoginBar loppuun, Footer, Loading, Authenticate  scss, mixinien luontia ja kokeilua"

# --- Add a new day-log row 112 ---
$ws.Range("A110").Copy()
$ws.Range("A112").PasteSpecial(-4122)
$ws.Range("A112").Value = 44591

$ws.Range("B111").Copy()
$ws.Range("B112").PasteSpecial(-4122)
$ws.Range("B112").Value = 6

$ws.Range("C112").Value = "Buttons, DataView, Home, Inputs, Market, NavBar scss transitio ja muokkaus, NavBar refaktor, seuraa aktiivista sivua"
$ws.Range("D112").Value = "client"

# --- Extend the totals formulas to cover the new row ---
$ws.Range("B115").Formula = "=SUM(B2:B112)"

# --- Update the selection to reflect where the user ended up ---
[void]$ws.Range("D112").Select()

Write-Host "done"
